$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        if ($val -eq "System, dnasr281@gmail.com") {
            $cell.Value = "dnasr281@gmail.com, System"
        } elseif ($val -eq "dnasr281@gmail.com, System") {
            $cell.Value = "System, dnasr281@gmail.com"
        }
    }
}
